$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the status check timestamp in F1
$ws.Range("F1").Value = "Last status check on: 25.02.2022 16:30"

# Update row 7 (MOL Olomoucka) price data
$ws.Range("B7").Value = 39.5
$ws.Range("C7").Value = 38.5

# Force text type (leading apostrophe prevents numeric auto-conversion),
# then strip the auto-applied "Text" number format so the cell keeps the
# workbook's default style.
$ws.Range("D7").Value = "'+1.0"
$ws.Range("D7").ClearFormats()

$ws.Range("E7").Value = "2022-02-25 16:32:22"
$ws.Range("E7").ClearFormats()
